$wb = $excel.ActiveWorkbook

# Rename Sheet1 to AddCustomerTest
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "AddCustomerTest"

# Populate header row (A1:C1)
$ws1.Range("A1").Value = "firstname"
$ws1.Range("B1").Value = "lastname"
$ws1.Range("C1").Value = "postcode"

# Row2/Row3 first two columns
$ws1.Range("A2").Value = "Tai"
$ws1.Range("B2").Value = "Le"
$ws1.Range("A3").Value = "Giang"
$ws1.Range("B3").Value = "Nguyen"

# postcode column (C3 then C2)
$ws1.Range("C3").Value = "Abcd20"
$ws1.Range("C2").Value = "Test1000"

# alerttext column
$ws1.Range("D1").Value = "alerttext"
$ws1.Range("D2").Value = "Customer added successfully"
$ws1.Range("D3").Value = "Customer added successfully"

# Set selection to D3 as final active cell on sheet1
$ws1.Range("D3").Select()
